# The commit reshuffles the per-row facts (date, variety, quality, volume,
# min/max/avg price, origin, $/Kg) across rows 2-15 of the sheet, while the
# constant columns (A,B,C,E,F,G,N,Q,R) stay untouched. Capture the "before"
# state for every row first, then write it back out in the new order, so
# that writes never clobber a value we still need to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "H", "I", "J", "K", "L", "M", "O", "P")

# Snapshot the current (pre-edit) values for rows 2..15, column by column.
$before = @{}
foreach ($col in $cols) {
    $rowVals = @{}
    for ($r = 2; $r -le 15; $r++) {
        $rowVals[$r] = $ws.Range("$col$r").Value2
    }
    $before[$col] = $rowVals
}

# after-row -> before-row that supplies its new data.
$mapping = @{
    2  = 12
    3  = 7
    4  = 5
    5  = 6
    6  = 11
    7  = 9
    8  = 14
    9  = 15
    10 = 2
    11 = 13
    12 = 3
    13 = 4
    14 = 8
    15 = 10
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $before[$col][$srcRow]
    }
}
